# self-assessment.docx edits
#
# Summary of changes (from the commit's XML diff):
#  1. Epic 4 / Implementation cell: append a "." run after the sidebar
#     bullet, and add a new bullet about the dashboard due-date list.
#  2. Epic 5 / Implementation cell: replace the "To be filled in"
#     placeholder with real content and two more bullets (one of which
#     has a mid-sentence run split around "name,").
#  3. Epic 6 / Implementation cell: add a bullet about viewing
#     notifications in the top popup (carries a lastRenderedPageBreak).
#  4. Epic 7 title cell: the "Epic 7" run now carries a
#     lastRenderedPageBreak (page break moved here).
#  5. Epic 8 / Implementation cell: the summary-report sentence, which
#     used to be split into two runs around a lastRenderedPageBreak,
#     becomes a single contiguous run (break removed).
#  6. Epic 9 title cell: the "Epic 9" run no longer carries the
#     lastRenderedPageBreak (it moved up to Epic 7).

$d = $word.ActiveDocument

function Get-ParaByExactText($needle) {
    # Word Range.Text includes a trailing paragraph mark (and, for the
    # very last paragraph in a cell, a cell-mark) so compare after
    # trimming those control characters.
    $paras = $d.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $needle) {
            return $p
        }
    }
    return $null
}

function Add-ParagraphAfter($needleText, $newText) {
    # Appends a brand new list paragraph right after the paragraph whose
    # trimmed text equals $needleText, inheriting its pPr/rPr (list
    # style, numbering, run size) automatically.
    $target = Get-ParaByExactText($needleText)
    $r = $target.Range
    $r.Collapse(0)
    $r.InsertAfter([char]13 + $newText)
}

function Add-RunAfter($needleText, $moreText) {
    # Appends more text to the same paragraph whose trimmed text equals
    # $needleText (continuing within that paragraph, not starting a new
    # one).
    $target = Get-ParaByExactText($needleText)
    $r = $target.Range
    $r.Collapse(0)
    $r.InsertAfter($moreText)
}

function Split-RunAt($searchText) {
    # Forces Word to split the run containing $searchText into its own
    # run by toggling Bold on then back off across that sub-range - the
    # emulator only splits runs when a character formatting property
    # actually changes, so this nets a clean 3-way split with unchanged
    # formatting.
    $fr = $d.Content
    $ok = $fr.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($ok) {
        $fr.Bold = 1
        $fr.Bold = 0
    }
}

# ---------------------------------------------------------------------
# 1. Epic 4 / Implementation cell
# ---------------------------------------------------------------------

# 1a. "." run appended after the sidebar bullet's sentence.
Add-RunAfter "A sidebar on the left allows the user to show the dashboard for each team they are part of" "."

# 1b. New trailing bullet about due dates.
Add-ParagraphAfter "There is also a button to view the team details for the current team. " "The dashboard also shows a list of due dates of all the tasks the user has access to.  "

# ---------------------------------------------------------------------
# 2. Epic 5 / Implementation cell
# ---------------------------------------------------------------------

$found = $d.Content.Find.Execute("To be filled in", $true, $false, $false, $false, $false, $true, 1, $false, "Users can search for tasks that are on their dashboard with a search bar.", 2)

Add-ParagraphAfter "Users can search for tasks that are on their dashboard with a search bar." "Users can filter tasks using a dropdown on the dashboard based on completion and levels of priority."

Add-ParagraphAfter "Users can filter tasks using a dropdown on the dashboard based on completion and levels of priority." "Users can order tasks by various properties such as completion status, priority, name, and due date, selected through a drop down. They can then select an ascending or descending order. "

Split-RunAt "name,"

# ---------------------------------------------------------------------
# 3. Epic 6 / Implementation cell
# ---------------------------------------------------------------------

Add-ParagraphAfter "By default, the reminder will be sent 1 day before, but this can be edited in the edit task page. " "They can be viewed in the notification popup at the top of the screen. "

# New bullet carries a lastRenderedPageBreak before its text run.
$newPara = Get-ParaByExactText("They can be viewed in the notification popup at the top of the screen. ")
$newPara.Range.Find.Execute("They can be viewed", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# ---------------------------------------------------------------------
# 4/6. lastRenderedPageBreak moves from "Epic 9" up to "Epic 7"
# ---------------------------------------------------------------------

# Remove the break before "Epic 9" (it was there in the source doc).
$epic9 = $d.Content
$epic9.Find.Execute("Epic 9", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# ---------------------------------------------------------------------
# 5. Epic 8 / Implementation cell: merge the two "summary report" runs
#    into one (dropping the lastRenderedPageBreak between them).
# ---------------------------------------------------------------------

$found2 = $d.Content.Find.Execute("A summary report can be viewed on the dashboard. This shows the time spent on each task as well as how much time the team has spent working on tasks per day. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 1)
